$d = $word.ActiveDocument

# The document ends with two empty paragraphs right before the final
# section break. The second (last) one gains a run of text "paste"
# that starts with a lastRenderedPageBreak marker (as Word records when
# that run happens to start a new page).
$count = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($count)

# Collapse to the very start of that (empty) paragraph and graft in a
# minimal WordprocessingML fragment containing the run, so the existing
# paragraph mark/properties are left untouched.
$insertionPoint = $d.Range($p.Range.Start, $p.Range.Start)

$fragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:r>' +
                '<w:lastRenderedPageBreak/>' +
                '<w:t>paste</w:t>' +
              '</w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$insertionPoint.InsertXML($fragment)
